## Marks additional "Strings" / "Searching & Sorting" rows as completed ("yes")
## in column C, matching each row's existing Good/Neutral/Bad status-style
## (mirrors the preexisting conditional-style cells already on the sheet),
## and updates the current selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of target cell -> donor cell that already carries the desired
# "Good" (green), "Neutral" (yellow) or "Bad" (red) cell style, so the
# existing style indices in the workbook get reused instead of new ones
# being minted.
$targets = [ordered]@{
    "C59"  = "C6"    # Good
    "C87"  = "C15"   # Bad
    "C88"  = "C14"   # Neutral
    "C89"  = "C14"   # Neutral
    "C93"  = "C6"    # Good
    "C96"  = "C14"   # Neutral
    "C97"  = "C6"    # Good
    "C102" = "C6"    # Good
    "C103" = "C6"    # Good
    "C104" = "C6"    # Good
    "C105" = "C6"    # Good
    "C107" = "C6"    # Good
    "C108" = "C6"    # Good
    "C109" = "C6"    # Good
    "C110" = "C6"    # Good
    "C111" = "C6"    # Good
    "C112" = "C6"    # Good
    "C113" = "C6"    # Good
    "C114" = "C6"    # Good
    "C115" = "C6"    # Good
    "C116" = "C14"   # Neutral
}

foreach ($target in $targets.Keys) {
    $donor = $targets[$target]
    $ws.Range($donor).Copy()
    $ws.Range($target).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Range($target).Value = "yes"
}

# Reflect the place the editing session ended up at.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 355
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C116").Select()
